# Generate Report for Handback
# Refresh the handoff/handback timestamps recorded on the report sheets.

$wb = $excel.ActiveWorkbook

# "Overview" sheet: G3 = Latest HO Xliff Generate Date for the
# 68f3d6fe-... file (row in sync with en-US).
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G3").Value = "2016-09-02 14:56:26"

# "zh-cn" sheet: row 3 corresponds to the same 68f3d6fe-... file.
# H3 = Correspond Handoff Datetime, K3 = Correspond Handback DateTime.
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H3").Value = "2016-09-02 14:56:22"
$wsZhCn.Range("K3").Value = "2016-09-02 14:56:43"

# "de-de" sheet: row 3 corresponds to the same 68f3d6fe-... file.
# H3 = Correspond Handoff Datetime, K3 = Correspond Handback DateTime.
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H3").Value = "2016-09-02 14:56:26"
$wsDeDe.Range("K3").Value = "2016-09-02 14:56:51"
